$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new cells to be treated as plain text so that numeric-looking
# values (band_no/user_no) and the date string are not reinterpreted by
# Excel as numbers / date serials.
$ws.Range("A2:H3").NumberFormat = "@"

$ws.Range("A2").Value = "74311420"
$ws.Range("B2").Value = "80725555"
$ws.Range("C2").Value = "US"
$ws.Range("D2").Value = "2019-05-20"
$ws.Range("E2").Value = "invitation"
$ws.Range("F2").Value = "bandapp"
$ws.Range("G2").Value = "ios"
$ws.Range("H2").Value = "74311420"

$ws.Range("A3").Value = "74311420"
$ws.Range("B3").Value = "80763028"
$ws.Range("C3").Value = "US"
$ws.Range("D3").Value = "2019-05-22"
$ws.Range("E3").Value = "invitation"
$ws.Range("F3").Value = "bandapp"
$ws.Range("G3").Value = "ios"
$ws.Range("H3").Value = "74311420"
